$d = $word.ActiveDocument

# 1. Insert a new run containing a single leading space before the title
#    run. Using InsertXML (rather than InsertBefore, which would just grow
#    the existing run's text) forces a distinct <w:r> to be created, which
#    is what the target diff shows.
$titlePara = $d.Paragraphs(1)
$titleRange = $titlePara.Range
$insertPoint = $titleRange.Duplicate
$insertPoint.Collapse(1)
$spaceRunXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint.InsertXML($spaceRunXml)

# 2. Mark every run that carries a drawing (inline picture) as NoProofing,
#    which serializes as <w:rPr><w:noProof/></w:rPr> on that run.
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $range = $p.Range
    if ($range.InlineShapes.Count -gt 0) {
        $range.NoProofing = 1
    }
}
